$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "2023-06-26 19:33:36"

# Row 3
$ws.Range("A3").Value = "2023-06-26 19:33:36"

# Row 4
$ws.Range("A4").Value = "2023-06-26 19:33:36"

# Row 5
$ws.Range("A5").Value = "2023-06-26 19:33:36"
$ws.Range("D5").Value = 90306
$ws.Range("E5").Value = 67936
$ws.Range("F5").Value = 10885
$ws.Range("G5").Value = 3978
$ws.Range("H5").Value = 1519
$ws.Range("I5").Value = 5988

# Row 6
$ws.Range("A6").Value = "2023-06-26 19:33:37"
$ws.Range("D6").Value = 996382
$ws.Range("E6").Value = 770094
$ws.Range("F6").Value = 153163
$ws.Range("G6").Value = 29864
$ws.Range("H6").Value = 7386
$ws.Range("I6").Value = 35875

# Row 7
$ws.Range("A7").Value = "2023-06-26 19:33:37"
$ws.Range("D7").Value = 190322
$ws.Range("E7").Value = 128793
$ws.Range("F7").Value = 18885
$ws.Range("G7").Value = 8676
$ws.Range("H7").Value = 6061
$ws.Range("I7").Value = 27907

# Row 8
$ws.Range("A8").Value = "2023-06-26 19:33:37"
$ws.Range("D8").Value = 59628
$ws.Range("F8").Value = 4683
$ws.Range("G8").Value = 2482
$ws.Range("H8").Value = 2806
$ws.Range("I8").Value = 15768

# Row 9
$ws.Range("A9").Value = "2023-06-26 19:33:37"
$ws.Range("D9").Value = 774579
$ws.Range("E9").Value = 426307
$ws.Range("F9").Value = 71238
$ws.Range("G9").Value = 42261
$ws.Range("I9").Value = 203592

# Row 10
$ws.Range("A10").Value = "2023-06-26 19:33:38"
$ws.Range("D10").Value = 12599
$ws.Range("E10").Value = 5558
$ws.Range("F10").Value = 895
$ws.Range("G10").Value = 844
$ws.Range("H10").Value = 854
$ws.Range("I10").Value = 4448

# Row 11
$ws.Range("A11").Value = "2023-06-26 19:33:38"

# Row 12
$ws.Range("A12").Value = "2023-06-26 19:33:38"
$ws.Range("D12").Value = 635843
$ws.Range("E12").Value = 481462
$ws.Range("F12").Value = 106587
$ws.Range("G12").Value = 29969
$ws.Range("H12").Value = 4890
$ws.Range("I12").Value = 12935

# Row 13
$ws.Range("A13").Value = "2023-06-26 19:33:39"
$ws.Range("E13").Value = 92959
$ws.Range("F13").Value = 19819
$ws.Range("G13").Value = 7842
$ws.Range("H13").Value = 3940
$ws.Range("I13").Value = 18922

# Row 14
$ws.Range("A14").Value = "2023-06-26 19:33:39"

# Row 15
$ws.Range("A15").Value = "2023-06-26 19:33:39"
$ws.Range("D15").Value = 55353
$ws.Range("E15").Value = 34531
$ws.Range("F15").Value = 7787
$ws.Range("I15").Value = 8278

# Row 16
$ws.Range("A16").Value = "2023-06-26 19:33:40"
$ws.Range("D16").Value = 29614
$ws.Range("E16").Value = 19414
$ws.Range("F16").Value = 4317
$ws.Range("G16").Value = 1746
$ws.Range("I16").Value = 3195

# Row 17
$ws.Range("A17").Value = "2023-06-26 19:33:40"

# Row 18
$ws.Range("A18").Value = "2023-06-26 19:33:40"
$ws.Range("D18").Value = 189232
$ws.Range("E18").Value = 145886
$ws.Range("F18").Value = 28301
$ws.Range("G18").Value = 7666
$ws.Range("H18").Value = 1742
$ws.Range("I18").Value = 5637

# Row 19
$ws.Range("A19").Value = "2023-06-26 19:33:41"
$ws.Range("D19").Value = 277508
$ws.Range("E19").Value = 106057
$ws.Range("F19").Value = 20890
$ws.Range("G19").Value = 16280
$ws.Range("H19").Value = 19024
$ws.Range("I19").Value = 115257

# Row 20
$ws.Range("A20").Value = "2023-06-26 19:33:41"

# Row 21
$ws.Range("A21").Value = "2023-06-26 19:33:41"

# Row 22
$ws.Range("A22").Value = "2023-06-26 19:33:41"
$ws.Range("D22").Value = 379630
$ws.Range("E22").Value = 247619
$ws.Range("F22").Value = 36002
$ws.Range("G22").Value = 20635
$ws.Range("H22").Value = 12252
$ws.Range("I22").Value = 63122

# Row 23
$ws.Range("A23").Value = "2023-06-26 19:33:42"
$ws.Range("D23").Value = 36336
$ws.Range("E23").Value = 24114
$ws.Range("F23").Value = 4725
$ws.Range("G23").Value = 1951
$ws.Range("I23").Value = 4606

# Row 24
$ws.Range("A24").Value = "2023-06-26 19:33:42"
$ws.Range("D24").Value = 153507
$ws.Range("E24").Value = 80750
$ws.Range("F24").Value = 15514
$ws.Range("G24").Value = 9181
$ws.Range("H24").Value = 7712
$ws.Range("I24").Value = 40350

# Row 25
$ws.Range("A25").Value = "2023-06-26 19:33:42"

# Row 26
$ws.Range("A26").Value = "2023-06-26 19:33:43"

# Row 27
$ws.Range("A27").Value = "2023-06-26 19:33:43"

# Row 28
$ws.Range("A28").Value = "2023-06-26 19:33:43"

# Row 29
$ws.Range("A29").Value = "2023-06-26 19:33:43"

# Row 30
$ws.Range("A30").Value = "2023-06-26 19:33:44"

# Row 31
$ws.Range("A31").Value = "2023-06-26 19:33:44"

# Row 32
$ws.Range("A32").Value = "2023-06-26 19:33:44"

# Row 33
$ws.Range("A33").Value = "2023-06-26 19:33:45"

# Row 34
$ws.Range("A34").Value = "2023-06-26 19:33:45"

# Row 35
$ws.Range("A35").Value = "2023-06-26 19:33:45"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "3.7"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = 2455
$ws.Range("I35").Value = 562

# Row 36
$ws.Range("A36").Value = "2023-06-26 19:33:45"

# Row 37
$ws.Range("A37").Value = "2023-06-26 19:33:46"

# Row 38
$ws.Range("A38").Value = "2023-06-26 19:33:46"
